$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.353.81"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.03"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.79"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4709"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2880"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06454"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.17"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.19"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7253"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.130"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.83"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.344.83"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.00"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007502"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.113.10"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.242"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.233"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.67"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.054"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.69"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.877"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.321"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09627"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.484"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.217"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.112"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6877"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01883"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.810"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.221"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.16"
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4221"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.931"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8219"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.83"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.602"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.25"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.960"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "898.84"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05725"
